# Apply the "pruebas / pantalla telefonos" update to the "Pruebas" sheet:
#  - expand the project/screen titles
#  - capitalize the first letter of several description/result cells
#  - add a new test-case row (row 10) for the "Enter" key behaviour
#  - nudge the view/selection to match the saved workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pruebas")

# --- Header titles -------------------------------------------------------
$ws.Range("A1").Value = "Proyecto BlueWeb"
$ws.Range("A2").Value = "Pantalla Teléfonos"

# --- Row 4: Consulta datos grafica ---------------------------------------
$ws.Range("C4").Value = "Que la tabla represente los datos obtenidos de la base de datos de la tabla ciudades. "
$ws.Range("E4").Value = "Hace consulta a la base de datos  SELECT_H_ACTIVA_TIPO_TELEFONOS, donde se le envia la descripción, el clave y la telefonía. Trae un conteo de cada una de las ciudades)."
$ws.Range("F4").Value = "Grafica con datos mostrados con exito."

# --- Row 5: Boton eliminar ------------------------------------------------
$ws.Range("C5").Value = "Que la tabla elimine cada registro al momento de presionar el boton eliminar"
$ws.Range("F5").Value = "Elementos eliminados con exito "

# --- Row 6: Boton de agregar ----------------------------------------------
$ws.Range("C6").Value = "Que la tabla agregue un nuevo registro  la tabla "
$ws.Range("E6").Value = "Agrega registos a la tabla al llenar los campos solicitados"
$ws.Range("F6").Value = "Elementos agregados con éxito"

# --- Row 7: Boton de modificar ---------------------------------------------
$ws.Range("C7").Value = "Que cada registro se modifique en la tabla"
$ws.Range("E7").Value = "Modifica los campos de la tabla"
$ws.Range("F7").Value = "Elementos modificados exitosamente"

# --- Row 8: editar solo un dato de la fila ---------------------------------
$ws.Range("C8").Value = "Se intenta editar solo un campo"
$ws.Range("E8").Value = "Edita un campo y manda ese valor a la funcion de editar, ignorando los otros"
$ws.Range("F8").Value = "Edita un solo campo dejando los otros iguales"

# --- Row 9: Editar 2 datos de la fila ---------------------------------------
$ws.Range("C9").Value = "Se intenta editar 2 campos sin que cambie el tercero"
$ws.Range("E9").Value = "Edita 2 campos del registro dejando el tercero igual"
$ws.Range("F9").Value = "Edita solo los 2 campos exitosamente"

# --- Row 10: new test case (Presionar boton enter) --------------------------
# Copy the formatting of the row above (date column needs the m/d/yyyy
# number format / border that every other data row already carries) before
# filling in the values.
$ws.Range("A9:F9").Copy()
$ws.Range("A10:F10").PasteSpecial(-4122)

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Presionar boton enter"
$ws.Range("C10").Value = "Al presionar el boton enter al agregar se enviar el formulario correctamente"
$ws.Range("D10").Value = 44431
$ws.Range("E10").Value = "Cuando se abre el modal para eliminar, y se llenan lo campos, al presionar enter se activa el boton de agregar"
$ws.Range("F10").Value = "El boton de agregar, agrega correctamente al presionar enter"

# --- View state -------------------------------------------------------------
$ws.Activate()
$ws.Range("E8").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
